$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 11:16:03"
$wsZh.Range("H2").Value = "2016-03-23 11:16:34"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 11:16:14"
$wsDe.Range("H2").Value = "2016-03-23 11:16:46"
